$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Batting tables (rows 2-12): left block columns A-F, right block columns J-O ---
$ws.Range("B2").Value = 22
$ws.Range("C2").Value = 6
$ws.Range("D2").Value = "Bowled"
$ws.Range("E2").Value = " Mark Wood"
$ws.Range("K2").Value = 73
$ws.Range("L2").Value = 26
$ws.Range("M2").Value = "LBW"
$ws.Range("N2").Value = " Mitchell Starc"

$ws.Range("B3").Value = 9
$ws.Range("C3").Value = 6
$ws.Range("D3").Value = "LBW"
$ws.Range("E3").Value = " Chris Woakes"
$ws.Range("K3").Value = 4
$ws.Range("L3").Value = 5
$ws.Range("N3").Value = " Pat Cummins"

$ws.Range("B4").Value = 75
$ws.Range("C4").Value = 21
$ws.Range("D4").Value = "Bowled"
$ws.Range("E4").Value = " Mark Wood"
$ws.Range("K4").Value = 62
$ws.Range("L4").Value = 23

$ws.Range("B5").Value = 23
$ws.Range("C5").Value = 11
$ws.Range("E5").Value = " Chris Woakes"
$ws.Range("K5").Value = 7
$ws.Range("L5").Value = 4
$ws.Range("M5").Value = "LBW"
$ws.Range("N5").Value = " Pat Cummins"

$ws.Range("B6").Value = 50
$ws.Range("C6").Value = 18
$ws.Range("D6").Value = "Caught"
$ws.Range("K6").Value = 10
$ws.Range("L6").Value = 4
$ws.Range("N6").Value = " Marcus Stionis"

$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "LBW"
$ws.Range("E7").Value = " Chris Woakes"
$ws.Range("K7").Value = 13
$ws.Range("L7").Value = 6
$ws.Range("M7").Value = "Bowled"
$ws.Range("N7").Value = " Adam Zampa"

$ws.Range("B8").Value = 0
$ws.Range("C8").Value = 1
$ws.Range("E8").Value = " Chris Woakes"
$ws.Range("K8").Value = 52
$ws.Range("L8").Value = 16
$ws.Range("M8").Value = "* NOT OUT"
$ws.Range("N8").Value = " "

$ws.Range("B9").Value = 17
$ws.Range("C9").Value = 7
$ws.Range("D9").Value = "Caught"
$ws.Range("E9").Value = " Mark Wood"
$ws.Range("L9").Value = 4
$ws.Range("M9").Value = "NOT OUT"
$ws.Range("N9").Value = " "

$ws.Range("B10").Value = 2
$ws.Range("C10").Value = 3
$ws.Range("D10").Value = "LBW"
$ws.Range("E10").Value = " Mark Wood"
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").Value = " "

$ws.Range("B11").Value = 9
$ws.Range("C11").Value = 3
$ws.Range("D11").Value = "NOT OUT"
$ws.Range("E11").Value = " "

$ws.Range("B12").Value = 11
$ws.Range("C12").Value = 4

# --- Totals row 16 ---
$ws.Range("A16").Value = 223
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "13.4"
$ws.Range("C16").ClearFormats()
$ws.Range("D16").Value = 82
$ws.Range("J16").Value = 227
$ws.Range("K16").Value = 6
$ws.Range("L16").NumberFormat = "@"
$ws.Range("L16").Value = "14.4"
$ws.Range("L16").ClearFormats()
$ws.Range("M16").Value = 88

# --- Bowling table rows 21-25 ---
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "2.0"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").Value = 37
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 18.5
$ws.Range("J21").Value = "Josh Hazlewood"
$ws.Range("K21").NumberFormat = "@"
$ws.Range("K21").Value = "3.0"
$ws.Range("K21").ClearFormats()
$ws.Range("L21").Value = 50
$ws.Range("M21").Value = 1
$ws.Range("N21").Value = 16.67

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = "3.0"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").Value = 40
$ws.Range("D22").Value = 4
$ws.Range("E22").Value = 13.33
$ws.Range("J22").Value = "Adam Zampa"
$ws.Range("K22").NumberFormat = "@"
$ws.Range("K22").Value = "3.0"
$ws.Range("K22").ClearFormats()
$ws.Range("L22").Value = 47
$ws.Range("M22").Value = 1
$ws.Range("N22").Value = 15.67

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = "3.0"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").Value = 70
$ws.Range("D23").Value = 0
$ws.Range("E23").Value = 23.33
$ws.Range("J23").Value = "Mitchell Starc"
$ws.Range("K23").NumberFormat = "@"
$ws.Range("K23").Value = "3.0"
$ws.Range("K23").ClearFormats()
$ws.Range("L23").Value = 37
$ws.Range("M23").Value = 1
$ws.Range("N23").Value = 12.33

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = "3.0"
$ws.Range("B24").ClearFormats()
$ws.Range("C24").Value = 36
$ws.Range("D24").Value = 1
$ws.Range("E24").Value = 12
$ws.Range("J24").Value = "Pat Cummins"
$ws.Range("K24").NumberFormat = "@"
$ws.Range("K24").Value = "3.0"
$ws.Range("K24").ClearFormats()
$ws.Range("L24").Value = 38
$ws.Range("M24").Value = 2
$ws.Range("N24").Value = 12.67

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = "2.4"
$ws.Range("B25").ClearFormats()
$ws.Range("C25").Value = 40
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = 16.67
$ws.Range("J25").Value = "Marcus Stionis"
$ws.Range("K25").NumberFormat = "@"
$ws.Range("K25").Value = "2.4"
$ws.Range("K25").ClearFormats()
$ws.Range("L25").Value = 55
$ws.Range("M25").Value = 1
$ws.Range("N25").Value = 22.92
